$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly-logged progress report rows (20-26) with date, hours,
# and activity description. Setting .Value on column C with a new string
# appends it to the shared-strings table in the same order we write them,
# which reproduces the uniqueCount 21 -> 28 growth from the diff.

$rows = @(
    @{ Row = 20; Date = 42840; Hours = 3;  Activity = "Development - Laying framework for website design" },
    @{ Row = 21; Date = 42848; Hours = 6;  Activity = "Development - Getting started!" },
    @{ Row = 22; Date = 42849; Hours = 12; Activity = "Development - Creating pages" },
    @{ Row = 23; Date = 42850; Hours = 12; Activity = "Development - Creating pages, adding features" },
    @{ Row = 24; Date = 42851; Hours = 12; Activity = "Development - Beautification of website, fixing bugs" },
    @{ Row = 25; Date = 42852; Hours = 12; Activity = "Development - Writing Report, fixing pages" },
    @{ Row = 26; Date = 42853; Hours = 3;  Activity = "Presenting - Demonstration! And last minute fixes before deploying" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Range("B$rowNum").Value = $r.Hours
    $ws.Range("C$rowNum").Value = $r.Activity
}

# Copy the date-formatted style from an already-populated date cell (A19)
# onto the new date cells so they pick up the existing numFmtId=14 style
# (s="12") instead of creating a brand new custom style, then write in the
# serial date values.
$ws.Range("A19").Copy()
$ws.Range("A20:A26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

foreach ($r in $rows) {
    $ws.Range("A$($r.Row)").Value = $r.Date
}

# Select A1:C26 to match the updated sheet-view selection in the diff.
$ws.Range("A1:C26").Select()
